$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.971.97"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "1.760.67"
$ws.Range("E3").Value = "  -2.72%  "

$ws.Range("E4").Value = "  -0.71%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3782"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3355"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.63"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.121"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07188"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.191"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.197"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "1.759.63"
$ws.Range("E16").Value = "  -3.56%  "

$ws.Range("E17").Value = "  -4.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06573"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.37"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9993"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.98"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.280"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.46%  "

$ws.Range("D23").Value = "27.986.81"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.70"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.371"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.337"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.66%  "

$ws.Range("D29").Value = "1.960.34"
$ws.Range("E29").Value = "  -3.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.77"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.249"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -15.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.018"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.796"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08813"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.25"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02341"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6616"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.25%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.154"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.49%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06173"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2116"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.211"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.450"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -10.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.026"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.86"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6062"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.824"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.84"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("E49").Value = "  -6.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.185"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.62%  "

$ws.Range("E51").Value = "  -0.59%  "
